$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (dated 2022-11-11, serial 44876) is inserted
# right before the existing row 64, shifting every following Jengibre record
# (rows 64-77) down by one row (new rows 65-78). Insert a whole row at 64
# so everything below cascades down automatically, mirroring Excel's
# Rows(...).Insert behaviour.
$ws.Rows("64:64").Insert()

# Populate the freshly inserted row 64 with the new observation. All of the
# "constant" descriptive columns (market/region/category/etc.) match every
# other Jengibre row in this sheet.
$ws.Cells.Item(64, 1).Value = 8
$ws.Cells.Item(64, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(64, 3).Value = "Coquimbo"
$ws.Cells.Item(64, 4).Value = 44876
$ws.Cells.Item(64, 5).Value = 4
$ws.Cells.Item(64, 6).Value = 100114007
$ws.Cells.Item(64, 7).Value = "Jengibre"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 400
$ws.Cells.Item(64, 11).Value = 14000
$ws.Cells.Item(64, 12).Value = 15000
$ws.Cells.Item(64, 13).Value = 14500
$ws.Cells.Item(64, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(64, 15).Value = "Perú"
$ws.Cells.Item(64, 16).Value = 1115
$ws.Cells.Item(64, 17).Value = 13
$ws.Cells.Item(64, 18).Value = "Hortaliza"
